$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "WARN: replace failed for: $old"
    }
}

# ---------------------------------------------------------------------
# Change 1 (R4 name cell): merge "R" "4" "-" "Recomendación..." "." runs
# ---------------------------------------------------------------------
Replace-Text "R4-Recomendación de productos vendidos juntos." "R4-Recomendación de productos vendidos juntos."

# ---------------------------------------------------------------------
# Change 2 (R4 Resumen cell): split the summary sentence into three runs
# ("La solución debe de recomendar otros productos " / "basado en compras"
#  / " de artículos que normalmente son venidos juntos"), leaving the
# trailing "." run untouched.
# ---------------------------------------------------------------------
$rng = $d.Content.Duplicate
$found = $rng.Find.Execute("La solución debe de recomendar otros productos cuando el usuario haga una compra", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Text = "La solución debe de recomendar otros productos "

    $rng2 = $d.Range($rng.End, $rng.End)
    $rng2.InsertAfter("basado en compras")

    $rng3 = $d.Range($rng2.End, $rng2.End)
    $rng3.InsertAfter(" de artículos que normalmente son venidos juntos")

    # Forcing a tiny formatting toggle on each populated range keeps it as
    # its own run in the saved XML instead of being re-coalesced into its
    # neighbour (both have identical rPr otherwise).
    $rng.Font.Bold = 1
    $rng.Font.Bold = 0

    $rng2.Font.Bold = 1
    $rng2.Font.Bold = 0

    $rng3.Font.Bold = 1
    $rng3.Font.Bold = 0
} else {
    Write-Output "WARN: change2 find failed"
}

# ---------------------------------------------------------------------
# Change 3 (R5 name cell): merge "R" "5" "-" "Recomendación..." "." runs
# ---------------------------------------------------------------------
Replace-Text "R5-Recomendación de productos para clientes." "R5-Recomendación de productos para clientes."

# ---------------------------------------------------------------------
# Change 4 (R5 Resumen cell): merge the two summary runs + final "." run
# ---------------------------------------------------------------------
Replace-Text "El sistema debe de reportar que productos son más propensos a ser comprados por cierta empresa." "El sistema debe de reportar que productos son más propensos a ser comprados por cierta empresa."

# ---------------------------------------------------------------------
# Change 5 (R5 Entradas cell): merge "Datos de ventas " + rest + ". " runs
# ---------------------------------------------------------------------
Replace-Text "Datos de ventas de cada cliente a largo de toda la historia de la compañía. " "Datos de ventas de cada cliente a largo de toda la historia de la compañía. "

# ---------------------------------------------------------------------
# Change 6 (R5 Resultado cell): merge text + "." runs
# ---------------------------------------------------------------------
Replace-Text "Productos recomendados para cada cliente." "Productos recomendados para cada cliente."

# ---------------------------------------------------------------------
# Change 7 (R6 name cell): merge "R" "6" "-" "Clasificar..." into one run
# up to "diferentes ", then re-insert the _GoBack bookmark (moved from the
# tr level, see change 8) between that run and the final "categorías."
# run.
# ---------------------------------------------------------------------
$rng = $d.Content.Duplicate
$found = $rng.Find.Execute("R", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
$rng = $d.Content.Duplicate
$found = $rng.Find.Execute("Clasificar a los clientes en diferentes categorías", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    # Expand the range to also cover the leading "R" "6" "-" runs.
    $rng.MoveStart(1, -3)
    $rng.Text = "R6-Clasificar a los clientes en diferentes "

    $rng2 = $d.Range($rng.End, $rng.End)
    $rng2.InsertAfter("categorías")

    $rng.Font.Bold = 1
    $rng.Font.Bold = 0
    $rng2.Font.Bold = 1
    $rng2.Font.Bold = 0
} else {
    Write-Output "WARN: change7 find failed"
}

# ---------------------------------------------------------------------
# Change 8 (R6 Resumen cell): append " y que se adapte a su estilo de
# compra" as its own run right before the trailing "." run.
# ---------------------------------------------------------------------
$rng = $d.Content.Duplicate
$found = $rng.Find.Execute("El programa debe de clasificar a los clientes en diferentes categorías para así darles un servicio más personalizado", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $ins = $d.Range($rng.End, $rng.End)
    $ins.InsertAfter(" y que se adapte a su estilo de compra")
    $ins.Font.Bold = 1
    $ins.Font.Bold = 0
} else {
    Write-Output "WARN: change8 find failed"
}

# Move the _GoBack bookmark out of the row (it currently trails the
# Resumen cell at the <w:tr> level) so it only lives inside the Name
# cell (inserted above, in change 7).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$nameRng = $d.Content.Duplicate
$found = $nameRng.Find.Execute("R6-Clasificar a los clientes en diferentes ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $bm = $d.Range($nameRng.End, $nameRng.End)
    $d.Bookmarks.Add("_GoBack", $bm) | Out-Null
} else {
    Write-Output "WARN: bookmark reinsertion find failed"
}

# ---------------------------------------------------------------------
# Change 9 (R6 Entradas cell): merge "Da" + "tos de cada cliente..." runs
# ---------------------------------------------------------------------
Replace-Text "Datos de cada cliente durante toda su historia en " "Datos de cada cliente durante toda su historia en "

# ---------------------------------------------------------------------
# Change 10 (R6 Resultado cell): merge text + "." runs
# ---------------------------------------------------------------------
Replace-Text "Diferentes categorías de clientes." "Diferentes categorías de clientes."

Write-Output "done"
